# Scheduled runner update: refresh computed market-profit columns
# (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ) on the
# Anima_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with
# newly-pulled pricing data. Only cells H:N on the affected rows change.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 931634.9399999999
$ws.Range("I17").Value = 468.07144
$ws.Range("J17").Value = 1453088.4
$ws.Range("K17").Value = 1404.21432
$ws.Range("L17").Value = 4359265.199999999
$ws.Range("M17").Value = -1236.21432
$ws.Range("N17").Value = -4359601.199999999

$ws.Range("H33").Value = 322.75
$ws.Range("I33").Value = 342.76923
$ws.Range("K33").Value = 342.76923
$ws.Range("M33").Value = -113.76923

$ws.Range("H115").Value = 1590.125
$ws.Range("I115").Value = 1590.125
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 4770.375
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -3203.375

$ws.Range("H129").Value = 1560.1538
$ws.Range("J129").Value = 1983.6296
$ws.Range("L129").Value = 5950.8888
$ws.Range("N129").Value = -15950.8888

$ws.Range("H132").Value = 2290.9155
$ws.Range("I132").Value = 2293.3667
$ws.Range("J132").Value = 2277.5454
$ws.Range("K132").Value = 6880.1001
$ws.Range("L132").Value = 6832.6362
$ws.Range("M132").Value = -4350.1001
$ws.Range("N132").Value = -11892.6362

$ws.Range("H137").Value = 1578.5
$ws.Range("I137").Value = 1593.6957
$ws.Range("J137").Value = 1564.52
$ws.Range("K137").Value = 4781.0871
$ws.Range("L137").Value = 4693.559999999999
$ws.Range("M137").Value = -2231.0871
$ws.Range("N137").Value = -9793.559999999999

$ws.Range("H138").Value = 1574.67
$ws.Range("I138").Value = 629.9655
$ws.Range("J138").Value = 1960.5352
$ws.Range("K138").Value = 1889.8965
$ws.Range("L138").Value = 5881.6056
$ws.Range("M138").Value = 3250.1035
$ws.Range("N138").Value = -16161.6056

$ws.Range("H141").Value = 1682.9584
$ws.Range("I141").Value = 884.7174
$ws.Range("J141").Value = 20042.5
$ws.Range("K141").Value = 2654.1522
$ws.Range("L141").Value = 60127.5
$ws.Range("M141").Value = 2525.8478
$ws.Range("N141").Value = -70487.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1906
$ws.Range("I45").Value = 1677.9
$ws.Range("K45").Value = 1677.9
$ws.Range("M45").Value = -1300.9

$ws.Range("H61").Value = 5557358
$ws.Range("I61").Value = 6537378
$ws.Range("J61").Value = 3910
$ws.Range("K61").Value = 6537378
$ws.Range("L61").Value = 3910
$ws.Range("M61").Value = -6537166
$ws.Range("N61").Value = -4334

$ws.Range("H74").Value = 7247939
$ws.Range("I74").Value = 703.3953
$ws.Range("J74").Value = 19233752
$ws.Range("K74").Value = 703.3953
$ws.Range("L74").Value = 19233752
$ws.Range("M74").Value = 170.6047
$ws.Range("N74").Value = -19235500

$ws.Range("H77").Value = 7247939
$ws.Range("I77").Value = 703.3953
$ws.Range("J77").Value = 19233752
$ws.Range("K77").Value = 3516.9765
$ws.Range("L77").Value = 96168760
$ws.Range("M77").Value = 851.0234999999998
$ws.Range("N77").Value = -96177496

$ws.Range("H132").Value = 1028700.8
$ws.Range("I132").Value = 2501.7458
$ws.Range("J132").Value = 4812810
$ws.Range("K132").Value = 7505.2374
$ws.Range("L132").Value = 14438430
$ws.Range("M132").Value = -4975.2374
$ws.Range("N132").Value = -14443490

$ws.Range("H136").Value = 5557358
$ws.Range("I136").Value = 6537378
$ws.Range("J136").Value = 3910
$ws.Range("K136").Value = 19612134
$ws.Range("L136").Value = 11730
$ws.Range("M136").Value = -19609584
$ws.Range("N136").Value = -16830

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 44696.75
$ws.Range("I26").Value = 22577.428
$ws.Range("J26").Value = 75663.8
$ws.Range("K26").Value = 22577.428
$ws.Range("L26").Value = 75663.8
$ws.Range("M26").Value = -22285.428
$ws.Range("N26").Value = -76247.8

$ws.Range("H134").Value = 1878.2858
$ws.Range("I134").Value = 1873.0256
$ws.Range("J134").Value = 1898.8
$ws.Range("K134").Value = 5619.0768
$ws.Range("L134").Value = 5696.4
$ws.Range("M134").Value = -3084.0768
$ws.Range("N134").Value = -10766.4

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

$ws.Range("H138").Value = 50780
$ws.Range("J138").Value = 50780
$ws.Range("L138").Value = 50780
$ws.Range("N138").Value = -61060

$ws.Range("H140").Value = 80765
$ws.Range("J140").Value = 80765
$ws.Range("L140").Value = 80765
$ws.Range("N140").Value = -91125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 22444.215
$ws.Range("I23").Value = 6309.857
$ws.Range("J23").Value = 38578.57
$ws.Range("K23").Value = 6309.857
$ws.Range("L23").Value = 38578.57
$ws.Range("M23").Value = -6069.857
$ws.Range("N23").Value = -39058.57

$ws.Range("H27").Value = 22444.215
$ws.Range("I27").Value = 6309.857
$ws.Range("J27").Value = 38578.57
$ws.Range("K27").Value = 6309.857
$ws.Range("L27").Value = 38578.57
$ws.Range("M27").Value = -6117.857
$ws.Range("N27").Value = -38962.57

$ws.Range("H134").Value = 6761236.5
$ws.Range("I134").Value = 9620424
$ws.Range("K134").Value = 28861272
$ws.Range("M134").Value = -28858737

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 3538.5
$ws.Range("I105").Value = 176
$ws.Range("J105").Value = 4659.3335
$ws.Range("K105").Value = 528
$ws.Range("L105").Value = 13978.0005
$ws.Range("M105").Value = 2093
$ws.Range("N105").Value = -19220.0005

$ws.Range("H113").Value = 504.875
$ws.Range("I113").Value = 515.5
$ws.Range("J113").Value = 497.2857
$ws.Range("K113").Value = 1546.5
$ws.Range("L113").Value = 1491.8571
$ws.Range("M113").Value = 623.5
$ws.Range("N113").Value = -5831.8571

$ws.Range("H131").Value = 4008.8809
$ws.Range("I131").Value = 525.55554
$ws.Range("J131").Value = 4958.879
$ws.Range("K131").Value = 1576.66662
$ws.Range("L131").Value = 14876.637
$ws.Range("M131").Value = 3463.33338
$ws.Range("N131").Value = -24956.637

$ws.Range("H132").Value = 1691
$ws.Range("I132").Value = 942
$ws.Range("J132").Value = 2240.2666
$ws.Range("K132").Value = 8478
$ws.Range("L132").Value = 20162.3994
$ws.Range("M132").Value = -5948
$ws.Range("N132").Value = -25222.3994

$ws.Range("H136").Value = 3544
$ws.Range("I136").Value = 1182.5
$ws.Range("J136").Value = 6242.857
$ws.Range("K136").Value = 3547.5
$ws.Range("L136").Value = 18728.571
$ws.Range("M136").Value = 1552.5
$ws.Range("N136").Value = -28928.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("N39").Value = 0

$ws.Range("H132").Value = 1901.8088
$ws.Range("I132").Value = 1379.7037
$ws.Range("J132").Value = 3915.6428
$ws.Range("K132").Value = 4139.1111
$ws.Range("L132").Value = 11746.9284
$ws.Range("M132").Value = -1609.1111
$ws.Range("N132").Value = -16806.9284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2586.2856
$ws.Range("I82").Value = 2547.2942
$ws.Range("J82").Value = 2646.5454
$ws.Range("K82").Value = 2547.2942
$ws.Range("L82").Value = 2646.5454
$ws.Range("M82").Value = -2186.2942
$ws.Range("N82").Value = -3368.5454

$ws.Range("H85").Value = 2586.2856
$ws.Range("I85").Value = 2547.2942
$ws.Range("J85").Value = 2646.5454
$ws.Range("K85").Value = 2547.2942
$ws.Range("L85").Value = 2646.5454
$ws.Range("M85").Value = -1299.2942
$ws.Range("N85").Value = -5142.5454

$ws.Range("H106").Value = 95000
$ws.Range("J106").Value = 95000
$ws.Range("L106").Value = 95000
$ws.Range("N106").Value = -97524

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0

$ws.Range("H132").Value = 2796.6226
$ws.Range("I132").Value = 2540.2703
$ws.Range("K132").Value = 7620.8109
$ws.Range("M132").Value = -5090.8109

$ws.Range("H136").Value = 4275173.5
$ws.Range("I136").Value = 1525.375
$ws.Range("J136").Value = 11113010
$ws.Range("K136").Value = 4576.125
$ws.Range("L136").Value = 33339030
$ws.Range("M136").Value = -2026.125
$ws.Range("N136").Value = -33344130

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4558847
$ws.Range("I132").Value = 1522.3478
$ws.Range("J132").Value = 16205343
$ws.Range("K132").Value = 4567.0434
$ws.Range("L132").Value = 48616029
$ws.Range("M132").Value = -2037.0434
$ws.Range("N132").Value = -48621089

$ws.Range("H135").Value = 103699.72
$ws.Range("J135").Value = 103699.72
$ws.Range("L135").Value = 103699.72
$ws.Range("N135").Value = -113839.72

$ws.Range("H136").Value = 1252.07
$ws.Range("I136").Value = 1212.662
$ws.Range("J136").Value = 1348.5518
$ws.Range("K136").Value = 3637.986
$ws.Range("L136").Value = 4045.6554
$ws.Range("M136").Value = -1087.986
$ws.Range("N136").Value = -9145.6554

$ws.Range("H137").Value = 83571
$ws.Range("J137").Value = 83571
$ws.Range("L137").Value = 83571
$ws.Range("N137").Value = -93771

$ws.Range("H139").Value = 69632
$ws.Range("J139").Value = 73857.5
$ws.Range("L139").Value = 73857.5
$ws.Range("N139").Value = -84137.5
